# Moved figures to manuscript. Create statistics of papers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of task/paper statistics (row 70 intentionally left blank,
# consistent with the existing gaps between sections in the sheet).
$newRows = @(
    @{ Row = 71; Task = "Generate stats and figures on the number of published papers"; Pages = 1; Time = 10 },
    @{ Row = 72; Task = "Update the introduction and abstract"; Pages = 1; Time = 15 },
    @{ Row = 73; Task = "Update the methodology"; Pages = 1; Time = 48 },
    @{ Row = 74; Task = "Update section on improvement objective"; Pages = 1; Time = 97 },
    @{ Row = 75; Task = "Update the tables regarding functionality and technical problems"; Pages = 1; Time = 81 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Task

    $cellB = $ws.Cells.Item($r.Row, 2)
    $cellB.Value = $r.Pages
    $cellB.HorizontalAlignment = -4108

    $cellC = $ws.Cells.Item($r.Row, 3)
    $cellC.Value = $r.Time
    $cellC.HorizontalAlignment = -4108
}

# Update the view to reflect where the user ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 57
$ws.Range("D75").Select()
